# CHE_country_general.xlsx edit
# - Fix "Note" column placement for actual_emissions_elec_heat rows (188-217):
#   the Note value was stored one column too far right (M); move it to L.
# - Add new "daily_travel_time" parameter rows (218-249), one per year
#   1990-2021, with sourced values every 5-6 years (the "choke" rows) and
#   empty placeholder rows otherwise, plus a trailing blank spacer row 250.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Move the stray Note (column M) into column L for rows 188-217
# ---------------------------------------------------------------------
for ($r = 188; $r -le 217; $r++) {
    $src = $ws.Cells.Item($r, 13)   # M
    $dst = $ws.Cells.Item($r, 12)   # L
    $dst.Formula = $src.Formula
    $src.ClearContents()
}

# ---------------------------------------------------------------------
# 2) Append the daily_travel_time rows (218-249) + trailing blank row 250
# ---------------------------------------------------------------------
# Reference/format source rows already present in the sheet
$fmtG = $ws.Range("G188")   # Value column style (font sz 12)
$fmtK = $ws.Range("K188")   # Link column style (Hyperlink)

$unit = "hr"
$reference = "Federal Office for Spatial Development, Mobility and Transport Microcensus 2021"
$link = "https://www.are.admin.ch/are/de/home/mobilitaet/grundlagen-und-daten/mzmv.html"
$note = "Table G3.2.2.1, includes waiting times"

# Minutes-per-day figures sourced for specific years; rest are blank "choke" rows
$sourced = @{
    1994 = 82.6
    2000 = 93.3
    2005 = 97.5
    2010 = 91.7
    2015 = 90.4
    2021 = 80.2
}

$row = 218
for ($year = 1990; $year -le 2021; $year++) {
    $ws.Range("A$row").Value = "CHE"
    $ws.Range("B$row").Value = "country"
    $ws.Range("C$row").Value = "daily_travel_time"
    $ws.Range("D$row").Value = "annual"
    $ws.Range("E$row").Value = $year

    $fmtG.Copy()
    $ws.Range("G$row").PasteSpecial(-4122)
    $fmtK.Copy()
    $ws.Range("K$row").PasteSpecial(-4122)

    if ($sourced.ContainsKey($year)) {
        $minutes = $sourced[$year]
        $ws.Range("G$row").Formula = "=$minutes/60"
        $ws.Range("H$row").Value = $unit
        $ws.Range("J$row").Value = $reference
        $ws.Range("K$row").Value = $link
        $ws.Range("L$row").Value = $note
    }

    $ws.Rows.Item($row).RowHeight = 16

    $row = $row + 1
}

# Trailing blank spacer row (250): keeps the formatted G/K columns only
$fmtG.Copy()
$ws.Range("G250").PasteSpecial(-4122)
$fmtK.Copy()
$ws.Range("K250").PasteSpecial(-4122)
$ws.Rows.Item(250).RowHeight = 16

# ---------------------------------------------------------------------
# 3) Update the view: frozen-pane scroll position + active selection
# ---------------------------------------------------------------------
$ws.Range("C222").Select()

$excel.CutCopyMode = $false
